# ---------------------------------------------------------------------------
# Rebuild "B0CTTWZCVK_sales_po_comparison.xlsx" per the commit:
#   "Update with Correct Forecast output"
#
# 1. Rename Sheet1 -> "Sales vs PO", and add three more sheets:
#      "Weekly Growth", "Volume Insights", "Prediction Info"
# 2. "Sales vs PO": insert a new column C "Order Week" (= old ds column,
#    the PO-week date) while ds (col A) moves 6 days later per row, and the
#    PO_Requested_Qty column (now D) is zeroed out (moved to the new
#    "Weekly Growth" sheet).
# 3. "Weekly Growth": the previously-nonzero PO_Requested_Qty rows, each
#    with a week-over-week Growth% column.
# 4. "Volume Insights": summary stats (Total/Average/Max/Min PO quantity).
# 5. "Prediction Info": the predicted next-week PO quantity.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename existing sheet -----------------------------------------
$ws1 = $wb.ActiveSheet
$ws1.Name = "Sales vs PO"

# --- Add the three new sheets, each placed after the previous one ----------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Weekly Growth"

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Volume Insights"

$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "Prediction Info"

# --- Match the original sheet's page setup / outline conventions on the
# new sheets (outline summary rows below / columns to the right, and the
# 0.75in/1in/0.5in margins used throughout this workbook). --------------
foreach ($s in @($ws2, $ws3, $ws4)) {
    $s.Outline.SummaryRow = 1
    $s.Outline.SummaryColumn = 1
    $ps = $s.PageSetup
    $ps.LeftMargin = 54
    $ps.RightMargin = 54
    $ps.TopMargin = 72
    $ps.BottomMargin = 72
    $ps.HeaderMargin = 36
    $ps.FooterMargin = 36
}

# =============================================================================
# Sheet 1 ("Sales vs PO"): insert "Order Week" as column C, push
# PO_Requested_Qty to column D (all zeros), and shift the "ds" column (A)
# forward by 6 days for every data row.
# =============================================================================

$ws1.Cells.Item(1,1).Value = "ds"
$ws1.Cells.Item(1,2).Value = "y"
$ws1.Cells.Item(1,3).Value = "Order Week"
$ws1.Cells.Item(1,4).Value = "PO_Requested_Qty"
# New column D needs the same header styling (bold/centered/bordered) as
# the rest of row 1 -- copy formats across from an existing header cell.
$ws1.Range("A1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Columns: ds, y, Order Week (old ds), PO_Requested_Qty (now always 0)
$salesData = @(
    @(45333,0,45327,0),
    @(45340,0,45334,0),
    @(45347,1,45341,0),
    @(45354,2,45348,0),
    @(45361,3,45355,0),
    @(45368,4,45362,0),
    @(45375,0,45369,0),
    @(45382,0,45376,0),
    @(45389,2,45383,0),
    @(45396,2,45390,0),
    @(45403,0,45397,0),
    @(45410,4,45404,0),
    @(45417,0,45411,0),
    @(45424,1,45418,0),
    @(45431,1,45425,0),
    @(45438,1,45432,0),
    @(45445,1,45439,0),
    @(45452,4,45446,0),
    @(45459,1,45453,0),
    @(45466,1,45460,0),
    @(45473,2,45467,0),
    @(45480,4,45474,0),
    @(45487,21,45481,0),
    @(45494,7,45488,0),
    @(45501,1,45495,0),
    @(45508,2,45502,0),
    @(45515,6,45509,0),
    @(45522,1,45516,0),
    @(45529,16,45523,0),
    @(45536,5,45530,0),
    @(45543,1,45537,0),
    @(45550,3,45544,0),
    @(45557,8,45551,0),
    @(45564,2,45558,0),
    @(45571,1,45565,0),
    @(45578,5,45572,0),
    @(45585,1,45579,0),
    @(45599,1,45593,0),
    @(45606,4,45600,0),
    @(45613,3,45607,0),
    @(45620,0,45614,0),
    @(45627,0,45621,0),
    @(45634,3,45628,0),
    @(45641,5,45635,0),
    @(45648,3,45642,0),
    @(45655,3,45649,0)
)

$r = 2
foreach ($row in $salesData) {
    $ws1.Cells.Item($r,1).Value = $row[0]
    $ws1.Cells.Item($r,2).Value = $row[1]
    $ws1.Cells.Item($r,3).Value = $row[2]
    $ws1.Cells.Item($r,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws1.Cells.Item($r,4).Value = $row[3]
    $r = $r + 1
}

# =============================================================================
# Sheet 2 ("Weekly Growth"): the previously nonzero PO_Requested_Qty rows,
# with a week-over-week Growth% column.
# =============================================================================

$ws2.Cells.Item(1,1).Value = "ds"
$ws2.Cells.Item(1,2).Value = "PO_Requested_Qty"
$ws2.Cells.Item(1,3).Value = "Growth%"
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$growthData = @(
    @(45334,96,0),
    @(45341,24,-75),
    @(45348,24,0),
    @(45355,16,-33.33333333333334),
    @(45376,8,-50),
    @(45495,32,300),
    @(45509,16,-50),
    @(45537,16,0),
    @(45558,16,0)
)

$r = 2
foreach ($row in $growthData) {
    $ws2.Cells.Item($r,1).Value = $row[0]
    $ws2.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws2.Cells.Item($r,2).Value = $row[1]
    $ws2.Cells.Item($r,3).Value = $row[2]
    $r = $r + 1
}

# =============================================================================
# Sheet 3 ("Volume Insights"): summary stats.
# =============================================================================

$ws3.Cells.Item(1,1).Value = "Total_PO_Quantity"
$ws3.Cells.Item(1,2).Value = "Average_PO_Quantity"
$ws3.Cells.Item(1,3).Value = "Max_PO_Quantity"
$ws3.Cells.Item(1,4).Value = "Min_PO_Quantity"
$ws1.Range("A1:D1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Cells.Item(2,1).Value = 248
$ws3.Cells.Item(2,2).Value = 27.55555555555556
$ws3.Cells.Item(2,3).Value = 96
$ws3.Cells.Item(2,4).Value = 8

# =============================================================================
# Sheet 4 ("Prediction Info"): predicted next-week PO quantity.
# =============================================================================

$ws4.Cells.Item(1,1).Value = "Predicted_Next_Week_PO_Quantity"
$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws4.Cells.Item(2,1).Value = 0

# --- Select "Sales vs PO" as the active sheet to mirror the source file ----
$ws1.Activate()
